$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I and J - copy format (style) from H1, then set values
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill data rows 2-37: I = 1 (constant), J = copy of H (column 8)
for ($r = 2; $r -le 37; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
